$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '19.999.04'
$ws.Range('E2').Value = '  -4.66%  '

$ws.Range('D3').Value = '1.421.25'
$ws.Range('E3').Value = '  -4.91%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9988'
$ws.Range('E4').Value = '  -0.77%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9988'
$ws.Range('E5').Value = '  -0.67%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '276.42'
$ws.Range('E6').Value = '  -1.90%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3672'
$ws.Range('E7').Value = '  -3.36%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3103'
$ws.Range('E8').Value = '  -0.22%  '

$ws.Range('E9').Value = '  -6.62%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.048'
$ws.Range('E10').Value = '  +2.07%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06553'
$ws.Range('E11').Value = '  -4.63%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9989'
$ws.Range('E12').Value = '  -0.84%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.516'
$ws.Range('E13').Value = '  -0.01%  '

$ws.Range('E14').Value = '  -0.11%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.218'
$ws.Range('E15').Value = '  -1.69%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.419.48'
$ws.Range('E16').Value = '  -5.47%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001023'
$ws.Range('E17').Value = '  -3.39%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.05679'
$ws.Range('E18').Value = '  -13.20%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9984'
$ws.Range('E19').Value = '  -0.75%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.49'
$ws.Range('E20').Value = '  -12.34%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.632'
$ws.Range('E21').Value = '  -4.83%  '

$ws.Range('E22').Value = '  -1.13%  '

$ws.Range('E23').Value = '  +2.10%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.253'
$ws.Range('E24').Value = '  -4.00%  '

$ws.Range('D25').Value = '20.026.20'
$ws.Range('E25').Value = '  -4.54%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.285'
$ws.Range('E26').Value = '  -0.70%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '133.51'
$ws.Range('E27').Value = '  -8.97%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.36'
$ws.Range('E28').Value = '  -2.53%  '

$ws.Range('D29').Value = '1.577.71'
$ws.Range('E29').Value = '  -5.58%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '110.17'
$ws.Range('E30').Value = '  -2.82%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.924'
$ws.Range('E31').Value = '  -17.89%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.286'
$ws.Range('E32').Value = '  -8.71%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8247'
$ws.Range('E33').Value = '  -11.01%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07729'
$ws.Range('E34').Value = '  -2.06%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.482'
$ws.Range('E35').Value = '  +0.74%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '8.306'
$ws.Range('E36').Value = '  -1.87%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.946'
$ws.Range('E37').Value = '  -0.92%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05868'
$ws.Range('E38').Value = '  +2.21%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9976'
$ws.Range('E39').Value = '  -0.76%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02067'
$ws.Range('E40').Value = '  -1.60%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.52'
$ws.Range('E41').Value = '  -4.11%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1884'
$ws.Range('E42').Value = '  -3.65%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.103'
$ws.Range('E43').Value = '  -3.78%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.46'
$ws.Range('E44').Value = '  -2.12%  '

$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5333'
$ws.Range('E45').Value = '  -4.16%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.541'
$ws.Range('E46').Value = '  -2.73%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5206'
$ws.Range('E47').Value = '  -3.27%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '116.10'
$ws.Range('E48').Value = '  +3.31%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.778'
$ws.Range('E49').Value = '  -2.48%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.037'
$ws.Range('E50').Value = '  -7.64%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9983'
$ws.Range('E51').Value = '  -0.83%  '
